# Regenerate merged AHB files
# - Rename the "_old" / "_new" header suffixes to the specific format
#   version identifiers ("_FV2310" / "_FV2404").
# - Turn the data range into a real Excel Table (ListObject) so the sheet
#   gets an autofilter + structured table definition.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells -------------------------------------------------
$headerRenames = @{
    "A1" = "Segmentname_FV2310"
    "B1" = "Segmentgruppe_FV2310"
    "C1" = "Segment_FV2310"
    "D1" = "Datenelement_FV2310"
    "E1" = "Segment ID_FV2310"
    "F1" = "Code_FV2310"
    "G1" = "Qualifier_FV2310"
    "H1" = "Beschreibung_FV2310"
    "I1" = "Bedingungsausdruck_FV2310"
    "J1" = "Bedingung_FV2310"
    "L1" = "Segmentname_FV2404"
    "M1" = "Segmentgruppe_FV2404"
    "N1" = "Segment_FV2404"
    "O1" = "Datenelement_FV2404"
    "P1" = "Segment ID_FV2404"
    "Q1" = "Code_FV2404"
    "R1" = "Qualifier_FV2404"
    "S1" = "Beschreibung_FV2404"
    "T1" = "Bedingungsausdruck_FV2404"
    "U1" = "Bedingung_FV2404"
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value = $headerRenames[$addr]
}

# --- 2. Convert the used range A1:U68 into an Excel Table -------------------
$tableRange = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
